$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 900
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("H125").Value = 2310.5208
$ws.Range("I125").Value = 1402.909
$ws.Range("K125").Value = 12626.181
$ws.Range("M125").Value = -10166.181
$ws.Range("H127").Value = 1910.8182
$ws.Range("J127").Value = 1000
$ws.Range("L127").Value = 3000
$ws.Range("N127").Value = -12920
$ws.Range("H129").Value = 29412766
$ws.Range("I129").Value = 31250752
$ws.Range("J129").Value = 5000
$ws.Range("K129").Value = 93752256
$ws.Range("L129").Value = 15000
$ws.Range("M129").Value = -93747256
$ws.Range("N129").Value = -25000
$ws.Range("H131").Value = 4918.8667
$ws.Range("I131").Value = 3845.2727
$ws.Range("J131").Value = 7871.25
$ws.Range("K131").Value = 11535.8181
$ws.Range("L131").Value = 23613.75
$ws.Range("M131").Value = -6495.8181
$ws.Range("N131").Value = -33693.75
$ws.Range("H132").Value = 8495
$ws.Range("I132").Value = 9194.5
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 27583.5
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -25053.5
$ws.Range("N132").Value = -9560
$ws.Range("N32").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 3546.8
$ws.Range("I41").Value = 3546.8
$ws.Range("K41").Value = 3546.8
$ws.Range("M41").Value = -3132.8
$ws.Range("H132").Value = 1487.3478
$ws.Range("I132").Value = 888.65717
$ws.Range("K132").Value = 2665.97151
$ws.Range("M132").Value = -135.9715099999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 75395
$ws.Range("J132").Value = 75395
$ws.Range("L132").Value = 75395
$ws.Range("N132").Value = -85515
$ws.Range("H134").Value = 4049.7942
$ws.Range("J134").Value = 7619.9165
$ws.Range("L134").Value = 22859.7495
$ws.Range("N134").Value = -27929.7495
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H138").Value = 79489.75
$ws.Range("J138").Value = 79489.75
$ws.Range("L138").Value = 79489.75
$ws.Range("N138").Value = -89769.75
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("N141").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14534.798
$ws.Range("I31").Value = 1528.6364
$ws.Range("K31").Value = 1528.6364
$ws.Range("M31").Value = -1233.6364
$ws.Range("H34").Value = 14534.798
$ws.Range("I34").Value = 1528.6364
$ws.Range("K34").Value = 1528.6364
$ws.Range("M34").Value = -1326.6364
$ws.Range("H122").Value = 2326.5862
$ws.Range("I122").Value = 1782.591
$ws.Range("J122").Value = 4036.2856
$ws.Range("K122").Value = 5347.772999999999
$ws.Range("L122").Value = 12108.8568
$ws.Range("M122").Value = -2897.772999999999
$ws.Range("N122").Value = -17008.8568
$ws.Range("H132").Value = 44906.383
$ws.Range("I132").Value = 2344.2778
$ws.Range("K132").Value = 7032.8334
$ws.Range("M132").Value = -4502.8334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5735168.5
$ws.Range("I4").Value = 6396314.5
$ws.Range("K4").Value = 19188943.5
$ws.Range("M4").Value = -19188831.5
$ws.Range("H38").Value = 64
$ws.Range("I38").Value = 31.181818
$ws.Range("J38").Value = 154.25
$ws.Range("K38").Value = 93.54545400000001
$ws.Range("L38").Value = 462.75
$ws.Range("M38").Value = 253.454546
$ws.Range("N38").Value = -1156.75
$ws.Range("H51").Value = 6374.375
$ws.Range("I51").Value = 331.66666
$ws.Range("K51").Value = 994.9999799999999
$ws.Range("M51").Value = -534.9999799999999
$ws.Range("H113").Value = 2948.7778
$ws.Range("I113").Value = 5783.3335
$ws.Range("J113").Value = 2138.9048
$ws.Range("K113").Value = 17350.0005
$ws.Range("L113").Value = 6416.714399999999
$ws.Range("M113").Value = -15180.0005
$ws.Range("N113").Value = -10756.7144
$ws.Range("H122").Value = 921.2
$ws.Range("I122").Value = 1006.8333
$ws.Range("J122").Value = 792.75
$ws.Range("K122").Value = 9061.4997
$ws.Range("L122").Value = 7134.75
$ws.Range("M122").Value = -6611.4997
$ws.Range("N122").Value = -12034.75
$ws.Range("H128").Value = 186797.4
$ws.Range("I128").Value = 186797.4
$ws.Range("K128").Value = 560392.2
$ws.Range("M128").Value = -555412.2
$ws.Range("H132").Value = 1818.5
$ws.Range("I132").Value = 1117.2222
$ws.Range("J132").Value = 2239.2666
$ws.Range("K132").Value = 10054.9998
$ws.Range("L132").Value = 20153.3994
$ws.Range("M132").Value = -7524.9998
$ws.Range("N132").Value = -25213.3994

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7196
$ws.Range("I43").Value = 3008.5
$ws.Range("J43").Value = 9987.666999999999
$ws.Range("K43").Value = 3008.5
$ws.Range("L43").Value = 9987.666999999999
$ws.Range("M43").Value = -2857.5
$ws.Range("N43").Value = -10289.667
$ws.Range("H46").Value = 11369.8
$ws.Range("I46").Value = 7671.2856
$ws.Range("J46").Value = 19999.666
$ws.Range("K46").Value = 7671.2856
$ws.Range("L46").Value = 19999.666
$ws.Range("M46").Value = -7515.2856
$ws.Range("N46").Value = -20311.666
$ws.Range("H51").Value = 85083.336
$ws.Range("J51").Value = 85083.336
$ws.Range("L51").Value = 85083.336
$ws.Range("N51").Value = -86101.336
$ws.Range("H70").Value = 12500
$ws.Range("I70").Value = 12500
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 12500
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -12230
$ws.Range("H73").Value = 12500
$ws.Range("I73").Value = 12500
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 12500
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -11564
$ws.Range("H99").Value = 9950
$ws.Range("I99").Value = 9950
$ws.Range("K99").Value = 9950
$ws.Range("M99").Value = -7704
$ws.Range("H132").Value = 3912.92
$ws.Range("I132").Value = 3121.923
$ws.Range("J132").Value = 4769.8335
$ws.Range("K132").Value = 9365.769
$ws.Range("L132").Value = 14309.5005
$ws.Range("M132").Value = -6835.769
$ws.Range("N132").Value = -19369.5005
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3581.4
$ws.Range("I7").Value = 2105
$ws.Range("J7").Value = 9487
$ws.Range("K7").Value = 2105
$ws.Range("L7").Value = 9487
$ws.Range("M7").Value = -1993
$ws.Range("N7").Value = -9711
$ws.Range("H122").Value = 8681.857
$ws.Range("I122").Value = 4995
$ws.Range("K122").Value = 14985
$ws.Range("M122").Value = -12535
$ws.Range("H126").Value = 3581.4
$ws.Range("I126").Value = 2105
$ws.Range("J126").Value = 9487
$ws.Range("K126").Value = 6315
$ws.Range("L126").Value = 28461
$ws.Range("M126").Value = -3845
$ws.Range("N126").Value = -33401
$ws.Range("H132").Value = 6545.5415
$ws.Range("I132").Value = 7142.4287
$ws.Range("K132").Value = 21427.2861
$ws.Range("M132").Value = -18897.2861

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H114").Value = 49988
$ws.Range("J114").Value = 49988
$ws.Range("L114").Value = 49988
$ws.Range("N114").Value = -58666
$ws.Range("H126").Value = 3912.6667
$ws.Range("I126").Value = 3763.5715
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 11290.7145
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -8820.7145
$ws.Range("N126").Value = -22940
$ws.Range("H132").Value = 65259.188
$ws.Range("I132").Value = 11607.417
$ws.Range("J132").Value = 226214.5
$ws.Range("K132").Value = 34822.251
$ws.Range("L132").Value = 678643.5
$ws.Range("M132").Value = -32292.251
$ws.Range("N132").Value = -683703.5
$ws.Range("N31").ClearContents()
